$wb = $excel.ActiveWorkbook

# ALC row 31
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 766.3333
$ws.Range("I31").Value = 766.3333
$ws.Range("K31").Value = 2298.9999
$ws.Range("M31").Value = -2068.9999

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2244.1738
$ws.Range("I33").Value = 2126.2354
$ws.Range("J33").Value = 2578.3333
$ws.Range("K33").Value = 2126.2354
$ws.Range("L33").Value = 2578.3333
$ws.Range("M33").Value = -1897.2354
$ws.Range("N33").Value = -3036.3333

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 224.95653
$ws.Range("I39").Value = 156.16667
$ws.Range("K39").Value = 468.50001
$ws.Range("M39").Value = -172.50001

# ALC row 63
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""

# ALC row 66
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1140.4117
$ws.Range("I135").Value = 1045.3077
$ws.Range("J135").Value = 1449.5
$ws.Range("K135").Value = 9407.7693
$ws.Range("L135").Value = 13045.5
$ws.Range("M135").Value = -6872.7693
$ws.Range("N135").Value = -18115.5

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 440946.75
$ws.Range("I137").Value = 771892.9
$ws.Range("K137").Value = 2315678.7
$ws.Range("M137").Value = -2313128.7

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2768.6465
$ws.Range("I138").Value = 1259.8334
$ws.Range("J138").Value = 3103.9382
$ws.Range("K138").Value = 3779.5002
$ws.Range("L138").Value = 9311.8146
$ws.Range("M138").Value = 1360.4998
$ws.Range("N138").Value = -19591.8146

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = -724

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2294.44
$ws.Range("I32").Value = 2294.44
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2294.44
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2007.44
$ws.Range("N32").Value = ""

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 38464628
$ws.Range("I45").Value = 41669680
$ws.Range("K45").Value = 41669680
$ws.Range("M45").Value = -41669303

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 25866034
$ws.Range("I61").Value = 21742034
$ws.Range("J61").Value = 41674704
$ws.Range("K61").Value = 21742034
$ws.Range("L61").Value = 41674704
$ws.Range("M61").Value = -21741822
$ws.Range("N61").Value = -41675128

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 25866034
$ws.Range("I136").Value = 21742034
$ws.Range("J136").Value = 41674704
$ws.Range("K136").Value = 65226102
$ws.Range("L136").Value = 125024112
$ws.Range("M136").Value = -65223552
$ws.Range("N136").Value = -125029212

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -730

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1098.4445
$ws.Range("I94").Value = 735.75
$ws.Range("K94").Value = 735.75
$ws.Range("M94").Value = -284.75

# BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 89993
$ws.Range("J126").Value = 89993
$ws.Range("L126").Value = 89993
$ws.Range("N126").Value = -99873

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1681682.9
$ws.Range("J31").Value = 3340355.5
$ws.Range("L31").Value = 3340355.5
$ws.Range("N31").Value = -3340945.5

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1681682.9
$ws.Range("J34").Value = 3340355.5
$ws.Range("L34").Value = 3340355.5
$ws.Range("N34").Value = -3340759.5

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1476.2963
$ws.Range("I122").Value = 1371.7894
$ws.Range("K122").Value = 4115.3682
$ws.Range("M122").Value = -1665.3682

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6289.579
$ws.Range("I132").Value = 1904.625
$ws.Range("K132").Value = 5713.875
$ws.Range("M132").Value = -3183.875

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10783389
$ws.Range("I4").Value = 9004287
$ws.Range("J4").Value = 14341594
$ws.Range("K4").Value = 27012861
$ws.Range("L4").Value = 43024782
$ws.Range("M4").Value = -27012749
$ws.Range("N4").Value = -43025006

# CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 69.75
$ws.Range("J40").Value = 148.25
$ws.Range("L40").Value = 593
$ws.Range("N40").Value = -731

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2254.182
$ws.Range("I68").Value = 1879.6
$ws.Range("J68").Value = 6000
$ws.Range("K68").Value = 5638.799999999999
$ws.Range("L68").Value = 18000
$ws.Range("M68").Value = -4827.799999999999
$ws.Range("N68").Value = -19622

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2254.182
$ws.Range("I71").Value = 1879.6
$ws.Range("J71").Value = 6000
$ws.Range("K71").Value = 16916.4
$ws.Range("L71").Value = 54000
$ws.Range("M71").Value = -12860.4
$ws.Range("N71").Value = -62112

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3339.9333
$ws.Range("I139").Value = 2135.5881
$ws.Range("J139").Value = 4914.846
$ws.Range("K139").Value = 6406.7643
$ws.Range("L139").Value = 14744.538
$ws.Range("M139").Value = -1266.7643
$ws.Range("N139").Value = -25024.538

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 285.0909
$ws.Range("I2").Value = 217.25
$ws.Range("J2").Value = 466
$ws.Range("K2").Value = 217.25
$ws.Range("L2").Value = 466
$ws.Range("M2").Value = -104.25
$ws.Range("N2").Value = -692

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 13999.5
$ws.Range("J136").Value = 13999.5
$ws.Range("L136").Value = 41998.5
$ws.Range("N136").Value = -47098.5

# LTW row 42
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = ""

# LTW row 49
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = ""

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1922.6666
$ws.Range("I82").Value = 1173
$ws.Range("J82").Value = 2195.2727
$ws.Range("K82").Value = 1173
$ws.Range("L82").Value = 2195.2727
$ws.Range("M82").Value = -812
$ws.Range("N82").Value = -2917.2727

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1922.6666
$ws.Range("I85").Value = 1173
$ws.Range("J85").Value = 2195.2727
$ws.Range("K85").Value = 1173
$ws.Range("L85").Value = 2195.2727
$ws.Range("M85").Value = 75
$ws.Range("N85").Value = -4691.2727

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3642.8948
$ws.Range("J132").Value = 2663
$ws.Range("L132").Value = 7989
$ws.Range("N132").Value = -13049

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 78914.94
$ws.Range("I136").Value = 16689.182
$ws.Range("J136").Value = 192995.5
$ws.Range("K136").Value = 50067.546
$ws.Range("L136").Value = 578986.5
$ws.Range("M136").Value = -47517.546
$ws.Range("N136").Value = -584086.5
